# Weekly price-sheet update: add this week's two new Papaya price rows
# ("Primera" / "Segunda") at the top of the Femacal de La Calera block,
# pushing the five previously-existing rows down by two positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the existing data block (old rows
# 91-95). Excel's own shift semantics push that block down to 93-97 and
# carry its formatting (incl. the date-style column D) along for free.
$ws.Rows("91:92").Insert()

# New row 91: Papaya "Primera", week of 2023-08-28 (serial 45166)
$ws.Range("A91").Value = 3
$ws.Range("B91").Value = "Femacal de La Calera"
$ws.Range("C91").Value = "Coquimbo"
$ws.Range("D91").Value = 45166
$ws.Range("E91").Value = 5
$ws.Range("F91").Value = "Fruta"
$ws.Range("G91").Value = 100108
$ws.Range("H91").Value = "Tropicales y subtropicales"
$ws.Range("I91").Value = 100108004
$ws.Range("J91").Value = "Papaya"
$ws.Range("K91").Value = "Cultivar IV Región"
$ws.Range("L91").Value = "Primera"
$ws.Range("M91").Value = 45
$ws.Range("N91").Value = 20000
$ws.Range("O91").Value = 20000
$ws.Range("P91").Value = 20000
$ws.Range("Q91").Value = "$/bandeja 10 kilos"
$ws.Range("R91").Value = "Provincia del Elquí"
$ws.Range("S91").Value = 2000
$ws.Range("T91").Value = 10

# New row 92: Papaya "Segunda", same week (serial 45166)
$ws.Range("A92").Value = 3
$ws.Range("B92").Value = "Femacal de La Calera"
$ws.Range("C92").Value = "Coquimbo"
$ws.Range("D92").Value = 45166
$ws.Range("E92").Value = 5
$ws.Range("F92").Value = "Fruta"
$ws.Range("G92").Value = 100108
$ws.Range("H92").Value = "Tropicales y subtropicales"
$ws.Range("I92").Value = 100108004
$ws.Range("J92").Value = "Papaya"
$ws.Range("K92").Value = "Cultivar IV Región"
$ws.Range("L92").Value = "Segunda"
$ws.Range("M92").Value = 50
$ws.Range("N92").Value = 17000
$ws.Range("O92").Value = 17000
$ws.Range("P92").Value = 17000
$ws.Range("Q92").Value = "$/bandeja 10 kilos"
$ws.Range("R92").Value = "Provincia del Elquí"
$ws.Range("S92").Value = 1700
$ws.Range("T92").Value = 10
